$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Roll the timesheet forward by one week (Nov 2-8 -> Nov 9-15) ---
$ws.Range("B5").Value = 44144
$ws.Range("C5").Value = 44145
$ws.Range("D5").Value = 44146
$ws.Range("E5").Value = 44147
$ws.Range("F5").Value = 44148
$ws.Range("G5").Value = 44149
$ws.Range("H5").Value = 44150

# --- Clear last week's "Team Meting" hours for Tue (C8) and Thu (E8); ---
# --- the Daily Total (I8) and Week Total (I10) formulas recompute automatically ---
$ws.Range("C8").ClearContents()
$ws.Range("E8").ClearContents()

# --- Columns C:H were widened slightly (previously grouped with B at 6.77734375) ---
$ws.Range("C1:H1").EntireColumn.ColumnWidth = 7.109375

# --- Active cell moved to G7 ---
$ws.Range("G7").Select() | Out-Null
